# Add the newly trained "bert-mini" model section to the scores sheet.
# Strategy: duplicate the last existing section (ridge_classifier, rows 70-91)
# down to rows 93-114 (leaving row 92 blank as a separator, matching the
# existing layout convention), then edit the duplicated cells so they hold
# the bert-mini figures. Also bump the "weighted avg f1-score" value that
# already lived in B91.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Existing cell edit: weighted avg f1-score for ridge_classifier moves
#    from 0.8845 to 0.876
# ---------------------------------------------------------------------
$ws.Range("B91").Value = 0.876

# ---------------------------------------------------------------------
# 2. Duplicate rows 70:91 (A:F) down to row 93, preserving both the
#    values/text and the cell formatting (styles, borders, fills...).
# ---------------------------------------------------------------------
$src = $ws.Range("A70:F91")
$dst = $ws.Range("A93")

$src.Copy()
$dst.PasteSpecial(-4163)   # xlPasteValues
$src.Copy()
$dst.PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. Overwrite the header row with the new model's name + training time.
# ---------------------------------------------------------------------
$ws.Range("A93").Value = "bert-mini"
$ws.Range("F93").Value = "~240min"

# ---------------------------------------------------------------------
# 4. "appearance" classification report (rows 96-103)
# ---------------------------------------------------------------------
$ws.Range("B96").Value = 0.99
$ws.Range("C96").Value = 1
$ws.Range("D96").Value = 0.99
$ws.Range("E96").Value = 166837

$ws.Range("B97").Value = 0.81
$ws.Range("C97").Value = 0.82
$ws.Range("D97").Value = 0.82
$ws.Range("E97").Value = 166838

$ws.Range("B98").Value = 0.83
$ws.Range("C98").Value = 0.81
$ws.Range("D98").Value = 0.82
$ws.Range("E98").Value = 166837

$ws.Range("D100").Value = 0.88
$ws.Range("E100").Value = 500512

$ws.Range("B101").Value = 0.88
$ws.Range("C101").Value = 0.88
$ws.Range("D101").Value = 0.88
$ws.Range("E101").Value = 500512

$ws.Range("B102").Value = 0.88
$ws.Range("C102").Value = 0.88
$ws.Range("D102").Value = 0.88
$ws.Range("E102").Value = 500512

$ws.Range("B103").Value = 0.88480000000000003

# ---------------------------------------------------------------------
# 5. "palate" classification report (rows 107-114)
# ---------------------------------------------------------------------
$ws.Range("B107").Value = 0.96
$ws.Range("C107").Value = 1
$ws.Range("D107").Value = 0.98
$ws.Range("E107").Value = 166838

$ws.Range("B108").Value = 0.86
$ws.Range("C108").Value = 0.73
$ws.Range("D108").Value = 0.79
$ws.Range("E108").Value = 166837

$ws.Range("B109").Value = 0.79
$ws.Range("C109").Value = 0.88
$ws.Range("D109").Value = 0.83
$ws.Range("E109").Value = 166837

$ws.Range("D111").Value = 0.87
$ws.Range("E111").Value = 500512

$ws.Range("B112").Value = 0.87
$ws.Range("C112").Value = 0.87
$ws.Range("D112").Value = 0.87
$ws.Range("E112").Value = 500512

$ws.Range("B113").Value = 0.87
$ws.Range("C113").Value = 0.87
$ws.Range("D113").Value = 0.87
$ws.Range("E113").Value = 500512

$ws.Range("B114").Value = 0.86829999999999996

# ---------------------------------------------------------------------
# 6. Bring the new merged header/divider ranges in line with the rest of
#    the sheet's layout.
# ---------------------------------------------------------------------
$ws.Range("A93:E93").Merge()
$ws.Range("A94:E94").Merge()
$ws.Range("C103:E103").Merge()
$ws.Range("A104:E104").Merge()
$ws.Range("A105:E105").Merge()
$ws.Range("C114:E114").Merge()

# ---------------------------------------------------------------------
# 7. Update the view state to match where the user ended up scrolled to
#    after adding the new section.
# ---------------------------------------------------------------------
$ws.Range("B119").Select()
$excel.ActiveWindow.ScrollRow = 86
$excel.ActiveWindow.ScrollColumn = 1
